$d = $word.ActiveDocument

$replacements = @(
    @("598÷4=", "934÷9="),
    @("561÷6=", "491÷2="),
    @("505÷8=", "257÷2="),
    @("587÷9=", "380÷6="),
    @("134÷3=", "751÷9="),
    @("267÷4=", "623÷8="),
    @("291÷6=", "515÷2="),
    @("932÷6=", "581÷4="),
    @("239÷9=", "299÷6="),
    @("900÷2=", "986÷7="),
    @("917÷4=", "941÷4="),
    @("794÷8=", "948÷9="),
    @("328÷9=", "122÷7="),
    @("604÷4=", "490÷3="),
    @("843÷7=", "816÷4="),
    @("589÷6=", "879÷3="),
    @("607÷8=", "249÷4="),
    @("919÷3=", "897÷7="),
    @("729÷2=", "477÷2="),
    @("727÷9=", "683÷7="),
    @("488÷8=", "354÷3="),
    @("918÷2=", "224÷5="),
    @("926÷3=", "721÷2="),
    @("577÷3=", "717÷9="),
    @("274÷7=", "541÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
